# Correcting Relevance Markers Appenzeller-Herzog (2019) - van Dis (2020)
# Updates row 3 metrics values on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 0.5555555555555556
$ws.Range("D3").Value = 0.5555555555555556
$ws.Range("E3").Value = 0.8888888888888888
$ws.Range("F3").Value = 1

$ws.Range("H3").Value = 0.3717791411042945
$ws.Range("I3").Value = 0.2388100623757625
$ws.Range("J3").Value = 0.4444444444444444
$ws.Range("K3").Value = 392.1111111111111

$ws.Range("Q3").Value = 1052
$ws.Range("R3").Value = 20
$ws.Range("S3").Value = 152
$ws.Range("T3").Value = 476
$ws.Range("U3").Value = 745
$ws.Range("V3").Value = 569
$ws.Range("W3").Value = 1601
$ws.Range("X3").Value = 1469
$ws.Range("Y3").Value = 1145
$ws.Range("Z3").Value = 876

$ws.Range("AG3").Value = 0.987662
$ws.Range("AH3").Value = 0.906231
$ws.Range("AI3").Value = 0.706354
$ws.Range("AJ3").Value = 0.540407
